$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column A (unlabeled column duplicating the GENE values),
# shifting columns B:F left to become A:E.
$ws.Range("A:A").Delete()
